$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped from
# 45179 (2023-09-10) to 45180 (2023-09-11) for every data row (2-369).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 369 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
